# Update the Tournaments table: rename the "Submerged" season theme to
# "Unearthed" everywhere it appears in the tournament long names and in
# the generated OJS file names (columns B, D, E of the TournamentList
# table, rows 2-18).
$wb = $excel.ActiveWorkbook
$tournaments = $wb.Worksheets.Item("Tournaments")
$tournamentRange = $tournaments.Range("B2:E18")
$tournamentRange.Replace("Submerged", "Unearthed")

# Restore the previously-selected cell on the Tournaments sheet before
# switching focus away from it.
$tournaments.Activate()
$tournaments.Range("B4").Select()

# Make the Summary sheet the active (displayed) sheet/tab again, keeping
# its previous selection.
$summary = $wb.Worksheets.Item("Summary")
$summary.Activate()
